# Update the "想去人数" (F) and "最低票价" (G) figures for both the
# "展览" and "全部类型" sheets, mirroring the upstream data refresh.
#
# Row -> new F value, new G value (only set if the diff called for it)
$updates = @(
    @{ Row = 2;  F = 3148 },
    @{ Row = 6;  F = 1779 },
    @{ Row = 10; F = 10 },
    @{ Row = 11; F = 1471 },
    @{ Row = 12; F = 17 },
    @{ Row = 13; F = 579 },
    @{ Row = 15; F = 90 },
    @{ Row = 18; F = 68 },
    @{ Row = 21; F = 99 },
    @{ Row = 22; F = 117 },
    @{ Row = 23; F = 3464 },
    @{ Row = 25; F = 303 },
    @{ Row = 27; F = 90 },
    @{ Row = 28; F = 25 },
    @{ Row = 30; F = 1220; G = 65 },
    @{ Row = 31; F = 129 }
)

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($u in $updates) {
        $ws.Range("F" + $u.Row).Value = $u.F
        if ($u.ContainsKey("G")) {
            $ws.Range("G" + $u.Row).Value = $u.G
        }
    }
}
